$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-13 18:48:23'
$ws.Range('G2').Value = '100 cm'
$ws.Range('I2').Value = '2.1 mm'
$ws.Range('E3').Value = '2026-02-13 18:48:25'
$ws.Range('I3').Value = '5.7 mm'
$ws.Range('E4').Value = '2026-02-13 18:48:27'
$ws.Range('H4').Value = '''75%'
$ws.Range('I4').Value = '5.4 mm'
$ws.Range('J4').Value = '995.0 hPa'
$ws.Range('O4').Value = '9.4 °C'
$ws.Range('E5').Value = '2026-02-13 18:48:30'
$ws.Range('I5').Value = '0.6 mm'
$ws.Range('E6').Value = '2026-02-13 18:48:32'
$ws.Range('H6').Value = '''75%'
$ws.Range('I6').Value = '4.6 mm'
$ws.Range('J6').Value = '995.0 hPa'
$ws.Range('E7').Value = '2026-02-13 18:48:35'
$ws.Range('H7').Value = '''72%'
$ws.Range('I7').Value = '18.6 mm'
$ws.Range('J7').Value = '995.2 hPa'
$ws.Range('N7').Value = '10.7 °C 18:29 TU'
$ws.Range('E8').Value = '2026-02-13 18:48:37'
$ws.Range('H8').Value = '''79%'
$ws.Range('I8').Value = '20.2 mm'
$ws.Range('J8').Value = '995.2 hPa'
$ws.Range('N8').Value = '6.9 °C 18:25 TU'
$ws.Range('O8').Value = '9.4 °C'
$ws.Range('E9').Value = '2026-02-13 18:48:39'
$ws.Range('I9').Value = '3.6 mm'
$ws.Range('L9').Value = '23.4 km/h - 340º 18:27 TU'
$ws.Range('O9').Value = '9.5 °C'
$ws.Range('E10').Value = '2026-02-13 18:48:42'
$ws.Range('I10').Value = '17.9 mm'
$ws.Range('E11').Value = '2026-02-13 18:48:44'
$ws.Range('H11').Value = '''92%'
$ws.Range('I11').Value = '15.2 mm'
$ws.Range('E12').Value = '2026-02-13 18:48:46'
$ws.Range('H12').Value = '''84%'
$ws.Range('I12').Value = '5.1 mm'
$ws.Range('E13').Value = '2026-02-13 18:48:49'
$ws.Range('H13').Value = '''92%'
$ws.Range('J13').Value = '998.0 hPa'
$ws.Range('E14').Value = '2026-02-13 18:48:51'
$ws.Range('I14').Value = '19.6 mm'
$ws.Range('O14').Value = '10.5 °C'
$ws.Range('E15').Value = '2026-02-13 18:48:54'
$ws.Range('H15').Value = '''76%'
$ws.Range('I15').Value = '3.1 mm'
$ws.Range('E16').Value = '2026-02-13 18:48:56'
$ws.Range('I16').Value = '12.7 mm'
$ws.Range('E17').Value = '2026-02-13 18:48:59'
$ws.Range('I17').Value = '4.8 mm'
$ws.Range('O17').Value = '0.5 °C'
$ws.Range('E18').Value = '2026-02-13 18:49:01'
$ws.Range('H18').Value = '''82%'
$ws.Range('I18').Value = '8.9 mm'
$ws.Range('J18').Value = '995.2 hPa'
$ws.Range('E19').Value = '2026-02-13 18:49:04'
$ws.Range('H19').Value = '''89%'
$ws.Range('I19').Value = '13.7 mm'
$ws.Range('E20').Value = '2026-02-13 18:49:06'
$ws.Range('I20').Value = '20.7 mm'
$ws.Range('E21').Value = '2026-02-13 18:49:08'
$ws.Range('J21').Value = '998.0 hPa'
$ws.Range('E22').Value = '2026-02-13 18:49:10'
$ws.Range('E23').Value = '2026-02-13 18:49:13'
$ws.Range('I23').Value = '9.2 mm'
$ws.Range('E24').Value = '2026-02-13 18:49:15'
$ws.Range('J24').Value = '995.8 hPa'
$ws.Range('E25').Value = '2026-02-13 18:49:18'
$ws.Range('G25').Value = '112 cm'
$ws.Range('I25').Value = '8.8 mm'
$ws.Range('L25').Value = '46.4 km/h - 259º 18:25 TU'
$ws.Range('E26').Value = '2026-02-13 18:49:20'
$ws.Range('E27').Value = '2026-02-13 18:49:23'
$ws.Range('E28').Value = '2026-02-13 18:49:25'
$ws.Range('J28').Value = '995.4 hPa'
$ws.Range('E29').Value = '2026-02-13 18:49:28'
$ws.Range('H29').Value = '''86%'
$ws.Range('E30').Value = '2026-02-13 18:49:30'
$ws.Range('I30').Value = '4.0 mm'
$ws.Range('J30').Value = '994.9 hPa'
$ws.Range('E31').Value = '2026-02-13 18:49:33'
$ws.Range('H31').Value = '''72%'
$ws.Range('I31').Value = '2.9 mm'
$ws.Range('J31').Value = '993.9 hPa'
$ws.Range('O31').Value = '10.4 °C'
$ws.Range('E32').Value = '2026-02-13 18:49:35'
$ws.Range('O32').Value = '5.1 °C'
$ws.Range('E33').Value = '2026-02-13 18:49:38'
$ws.Range('I33').Value = '4.9 mm'
$ws.Range('J33').Value = '996.9 hPa'
$ws.Range('E34').Value = '2026-02-13 18:49:40'
$ws.Range('G34').Value = '107 cm'
$ws.Range('I34').Value = '10.5 mm'
$ws.Range('L34').Value = '49.3 km/h - 17º 18:14 TU'
$ws.Range('E35').Value = '2026-02-13 18:49:43'
$ws.Range('H35').Value = '''75%'
$ws.Range('I35').Value = '7.8 mm'
$ws.Range('J35').Value = '995.8 hPa'
$ws.Range('L35').Value = '76.3 km/h - 272º 18:00 TU'
$ws.Range('N35').Value = '3.6 °C 18:22 TU'
$ws.Range('O35').Value = '6.1 °C'
$ws.Range('E36').Value = '2026-02-13 18:49:45'
$ws.Range('I36').Value = '8.1 mm'
$ws.Range('J36').Value = '995.1 hPa'
$ws.Range('E37').Value = '2026-02-13 18:49:48'
$ws.Range('J37').Value = '996.9 hPa'
$ws.Range('O37').Value = '3.7 °C'
$ws.Range('E38').Value = '2026-02-13 18:49:50'
$ws.Range('I38').Value = '12.7 mm'
$ws.Range('E39').Value = '2026-02-13 18:49:53'
$ws.Range('I39').Value = '18.6 mm'
$ws.Range('E40').Value = '2026-02-13 18:49:55'
$ws.Range('J40').Value = '998.4 hPa'
$ws.Range('O40').Value = '1.5 °C'
$ws.Range('E41').Value = '2026-02-13 18:49:58'
$ws.Range('J41').Value = '995.3 hPa'
$ws.Range('E42').Value = '2026-02-13 18:50:00'
$ws.Range('E43').Value = '2026-02-13 18:50:03'
$ws.Range('I43').Value = '13.0 mm'
$ws.Range('E44').Value = '2026-02-13 18:50:05'
$ws.Range('H44').Value = '''90%'
$ws.Range('I44').Value = '5.9 mm'
$ws.Range('E45').Value = '2026-02-13 18:50:08'
$ws.Range('J45').Value = '993.6 hPa'
$ws.Range('O45').Value = '5.8 °C'
$ws.Range('E46').Value = '2026-02-13 18:50:10'
$ws.Range('H46').Value = '''90%'
$ws.Range('J46').Value = '995.9 hPa'
$ws.Range('L46').Value = '28.8 km/h - 321º 18:13 TU'
$ws.Range('O46').Value = '8.8 °C'
